$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 958
$ws1.Range("F3").Value = 1877
$ws1.Range("F4").Value = 421

# Sheet "全部类型" - update "想去人数" (F column) values for same events
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 958
$ws4.Range("F5").Value = 1877
$ws4.Range("F6").Value = 421
